# classcodes.xlsx: add a new "E" column of class-code regexes (the
# "LPHnumeric-only" variant of the existing D column) to the charlson
# sheet, rows 2-18, and switch the active sheet/selection to charlson
# (matches commit "Add possibility to remove all non LPHnumeric
# characters from codes. Fix #53.").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "charlson" sheet (position 2 in the tab strip)

# E2's value is identical to D2's ("41[02]"), and in the target workbook it
# also carries D2's cell style (s="7", horizontal-left alignment). Copy the
# format over first, then overwrite the value below.
$ws.Range('D2').Copy()
$ws.Range('E2').PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# Rows whose E value reuses a string already present elsewhere in the
# workbook (so no new shared-string entries are minted for these).
$ws.Range('E2').Value = '41[02]'
$ws.Range('E9').Value = '53[1-4]'
$ws.Range('E17').Value = '19[6-9]'
$ws.Range('E18').Value = '04[2-4]'

# Remaining rows need brand-new shared strings. They must be written in
# this exact sequence so the new entries land at shared-string indices
# 404-416 in the same order as the target workbook.
$ws.Range('E3').Value = '39891|4(0(2([01]1|91)|4([019][13]))|2(5[4-9]|8))'
$ws.Range('E4').Value = '0930|4(373|4[01]|3[1-9]|471)|557[19]|V434'
$ws.Range('E5').Value = '36234|43[0-8]'
$ws.Range('E6').Value = '29(0|41)|3312'
$ws.Range('E7').Value = '4(16[89]|90)|50([0-5]|64|8[18])'
$ws.Range('E8').Value = '4465|7(1(0[0-4]|4[0-28])|25)'
$ws.Range('E10').Value = '070([23]{2}|[45]4|[69])|57([01]|3[3489])|V427'
$ws.Range('E11').Value = '250[0-389]'
$ws.Range('E14').Value = '250[4-7]'
$ws.Range('E12').Value = '3(341|4([23]|4[0-69]))'
$ws.Range('E13').Value = '40(3([019]1)|4([019][23]))|58([256]|3[0-7]|80)|V(4(20|51)|56)'
$ws.Range('E15').Value = '1([4-68]|7[0-24-9]|9([0-4]|5[0-8]))|2(0[0-8]|386)'
$ws.Range('E16').Value = '456[0-2]|572[2-8]'

# Make "charlson" the active sheet/tab with E21 selected (this also clears
# tabSelected/activeTab from whichever sheet was previously active, i.e.
# "hip_ae").
$ws.Activate()
$ws.Range('E21').Select()
